# Scheduled-runner update: refresh cached market-board pricing/profit
# figures (columns H:N) across several Leve-profit sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 261.4
$ws.Range("I41").Value = 301.75
$ws.Range("J41").Value = 100
$ws.Range("K41").Value = 301.75
$ws.Range("L41").Value = 100
$ws.Range("M41").Value = 138.25
$ws.Range("N41").Value = -980

$ws.Range("H53").Value = 104.22222
$ws.Range("I53").Value = 64.92308
$ws.Range("J53").Value = 206.4
$ws.Range("K53").Value = 64.92308
$ws.Range("L53").Value = 206.4
$ws.Range("M53").Value = 572.07692
$ws.Range("N53").Value = -1480.4

$ws.Range("H62").Value = 4191.0386
$ws.Range("I62").Value = 1546.6364
$ws.Range("J62").Value = 6130.2666
$ws.Range("K62").Value = 1546.6364
$ws.Range("L62").Value = 6130.2666
$ws.Range("M62").Value = -922.6364000000001
$ws.Range("N62").Value = -7378.2666

$ws.Range("H65").Value = 4191.0386
$ws.Range("I65").Value = 1546.6364
$ws.Range("J65").Value = 6130.2666
$ws.Range("K65").Value = 7733.182000000001
$ws.Range("L65").Value = 30651.333
$ws.Range("M65").Value = -4613.182000000001
$ws.Range("N65").Value = -36891.333

$ws.Range("H96").Value = 100001350
$ws.Range("I96").Value = 111112530
$ws.Range("J96").Value = 800
$ws.Range("K96").Value = 333337590
$ws.Range("L96").Value = 2400
$ws.Range("M96").Value = -333336217
$ws.Range("N96").Value = -5146

$ws.Range("H97").Value = 1559.875
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 1559.875
$ws.Range("K97").Value = 0
$ws.Range("L97").ClearContents()
$ws.Range("M97").Value = 4679.625
$ws.Range("N97").Value = -5671.625

$ws.Range("H107").Value = 1766.6
$ws.Range("I107").Value = 1766.6
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 1766.6
$ws.Range("L107").Value = 0
$ws.Range("M107").ClearContents()
$ws.Range("N107").Value = 153.4000000000001

$ws.Range("H112").Value = 3450
$ws.Range("I112").Value = 3750
$ws.Range("J112").Value = 3375
$ws.Range("K112").Value = 11250
$ws.Range("L112").Value = 10125
$ws.Range("M112").Value = -10142
$ws.Range("N112").Value = -12341

$ws.Range("H137").Value = 6262.967
$ws.Range("I137").Value = 4791.737
$ws.Range("J137").Value = 8804.182000000001
$ws.Range("K137").Value = 14375.211
$ws.Range("L137").Value = 26412.546
$ws.Range("M137").Value = -11825.211
$ws.Range("N137").Value = -31512.546

$ws.Range("H138").Value = 4478.2607
$ws.Range("I138").Value = 3500
$ws.Range("J138").Value = 6000
$ws.Range("K138").Value = 10500
$ws.Range("L138").Value = 18000
$ws.Range("M138").Value = -5360
$ws.Range("N138").Value = -28280

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H28").Value = 8503.700000000001
$ws.Range("I28").Value = 8503.700000000001
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 8503.700000000001
$ws.Range("L28").Value = 0
$ws.Range("M28").Value = -8311.700000000001

$ws.Range("H32").Value = 234856.11
$ws.Range("I32").Value = 1108.75
$ws.Range("J32").Value = 1436985.4
$ws.Range("K32").Value = 1108.75
$ws.Range("L32").Value = 1436985.4
$ws.Range("M32").Value = -821.75
$ws.Range("N32").Value = -1437559.4

$ws.Range("H74").Value = 1378.8485
$ws.Range("I74").Value = 1075.1072
$ws.Range("J74").Value = 3079.8
$ws.Range("K74").Value = 1075.1072
$ws.Range("L74").Value = 3079.8
$ws.Range("M74").Value = -201.1071999999999
$ws.Range("N74").Value = -4827.8

$ws.Range("H77").Value = 1378.8485
$ws.Range("I77").Value = 1075.1072
$ws.Range("J77").Value = 3079.8
$ws.Range("K77").Value = 5375.536
$ws.Range("L77").Value = 15399
$ws.Range("M77").Value = -1007.536
$ws.Range("N77").Value = -24135

$ws.Range("H99").Value = 8503.700000000001
$ws.Range("I99").Value = 8503.700000000001
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 8503.700000000001
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -5508.700000000001

$ws.Range("H101").Value = 40000
$ws.Range("I101").Value = 0
$ws.Range("J101").Value = 40000
$ws.Range("K101").Value = 0
$ws.Range("L101").Value = 40000
$ws.Range("N101").Value = -46490

$ws.Range("H122").Value = 628330.4399999999
$ws.Range("I122").Value = 717420.5
$ws.Range("J122").Value = 4700
$ws.Range("K122").Value = 2152261.5
$ws.Range("L122").Value = 14100
$ws.Range("M122").Value = -2149811.5
$ws.Range("N122").Value = -19000

$ws.Range("H132").Value = 2332.077
$ws.Range("I132").Value = 1320.8
$ws.Range("J132").Value = 5703
$ws.Range("K132").Value = 3962.4
$ws.Range("L132").Value = 17109
$ws.Range("M132").Value = -1432.4
$ws.Range("N132").Value = -22169

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1586.1111
$ws.Range("I86").Value = 1483.091
$ws.Range("J86").Value = 2039.4
$ws.Range("K86").Value = 1483.091
$ws.Range("L86").Value = 2039.4
$ws.Range("M86").Value = -360.0909999999999
$ws.Range("N86").Value = -4285.4

$ws.Range("H89").Value = 1586.1111
$ws.Range("I89").Value = 1483.091
$ws.Range("J89").Value = 2039.4
$ws.Range("K89").Value = 7415.455
$ws.Range("L89").Value = 10197
$ws.Range("M89").Value = -1799.455
$ws.Range("N89").Value = -21429

$ws.Range("H99").Value = 1179.05
$ws.Range("I99").Value = 954.55554
$ws.Range("J99").Value = 3199.5
$ws.Range("K99").Value = 954.55554
$ws.Range("L99").Value = 3199.5
$ws.Range("M99").Value = 543.44446
$ws.Range("N99").Value = -6195.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 211200
$ws.Range("I22").Value = 263250
$ws.Range("J22").Value = 3000
$ws.Range("K22").Value = 263250
$ws.Range("L22").Value = 3000
$ws.Range("M22").Value = -262900
$ws.Range("N22").Value = -3700

$ws.Range("H31").Value = 3238.2126
$ws.Range("I31").Value = 4859.9165
$ws.Range("J31").Value = 2682.2
$ws.Range("K31").Value = 4859.9165
$ws.Range("L31").Value = 2682.2
$ws.Range("M31").Value = -4564.9165
$ws.Range("N31").Value = -3272.2

$ws.Range("H34").Value = 3238.2126
$ws.Range("I34").Value = 4859.9165
$ws.Range("J34").Value = 2682.2
$ws.Range("K34").Value = 4859.9165
$ws.Range("L34").Value = 2682.2
$ws.Range("M34").Value = -4657.9165
$ws.Range("N34").Value = -3086.2

$ws.Range("H96").Value = 4064.4285
$ws.Range("I96").Value = 0
$ws.Range("J96").Value = 4064.4285
$ws.Range("K96").Value = 0
$ws.Range("L96").Value = 4064.4285
$ws.Range("N96").Value = -9556.4285

$ws.Range("H103").Value = 25759
$ws.Range("I103").Value = 16199
$ws.Range("J103").Value = 63999
$ws.Range("K103").Value = 16199
$ws.Range("L103").Value = 63999
$ws.Range("M103").Value = -15027
$ws.Range("N103").Value = -66343

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 0
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("L68").ClearContents()
$ws.Range("N68").Value = 0

$ws.Range("H71").Value = 0
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("L71").ClearContents()
$ws.Range("N71").Value = 0

$ws.Range("H114").Value = 662.8
$ws.Range("I114").Value = 751
$ws.Range("J114").Value = 530.5
$ws.Range("K114").Value = 2253
$ws.Range("L114").Value = 1591.5
$ws.Range("M114").Value = 1001
$ws.Range("N114").Value = -8099.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value = 500
$ws.Range("I20").Value = 500
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 500
$ws.Range("L20").ClearContents()
$ws.Range("M20").Value = -255
$ws.Range("N20").Value = 0

$ws.Range("H24").Value = 5115
$ws.Range("I24").Value = 5115
$ws.Range("J24").Value = 0
$ws.Range("K24").Value = 5115
$ws.Range("L24").Value = 0
$ws.Range("M24").ClearContents()
$ws.Range("N24").Value = -4942

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H51").Value = 31268.25
$ws.Range("I51").Value = 20000
$ws.Range("J51").Value = 35024.332
$ws.Range("K51").Value = 20000
$ws.Range("L51").Value = 35024.332
$ws.Range("M51").Value = -19490
$ws.Range("N51").Value = -36044.332

$ws.Range("H136").Value = 4013.3845
$ws.Range("I136").Value = 4245
$ws.Range("J136").Value = 1234
$ws.Range("K136").Value = 12735
$ws.Range("L136").Value = 3702
$ws.Range("M136").Value = -10185
$ws.Range("N136").Value = -8802
